$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOL")

# Refresh the test OrderId values used by the log-off / scroll test case.
# Leading apostrophe keeps these numeric-looking IDs stored as text, same
# as the existing OrderId values in this column.
$ws.Range("A2").Value = "'51499682"
$ws.Range("A3").Value = "'51492012"
$ws.Range("A4").Value = "'51492012"
